# Applies the edit described by the diff:
# A new data row was inserted at sheet row 98 (pushing the former rows
# 98..194 down to 99..195), and the new row 98 was populated with a new
# observation for "Haba" at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 98; this shifts rows 98:194 down to 99:195
# and Excel will extend the used range/dimension automatically.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record.
$ws.Cells.Item(98, 1).Value = 4
$ws.Cells.Item(98, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98, 3).Value = "Los Lagos"
$ws.Cells.Item(98, 4).Value = 45280
$ws.Cells.Item(98, 5).Value = 10
$ws.Cells.Item(98, 6).Value = 100112026
$ws.Cells.Item(98, 7).Value = "Haba"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 40
$ws.Cells.Item(98, 11).Value = 17000
$ws.Cells.Item(98, 12).Value = 17000
$ws.Cells.Item(98, 13).Value = 17000
$ws.Cells.Item(98, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(98, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(98, 16).Value = 680
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
